$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.296.46'
$ws.Range("E2").Value = '  +5.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.641.28'
$ws.Range("E3").Value = '  +5.27%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.76'
$ws.Range("E5").Value = '  +3.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.57'
$ws.Range("E6").Value = '  +3.67%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("E8").Value = '  +1.97%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.639.06'
$ws.Range("E9").Value = '  +5.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  +15.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.166'
$ws.Range("E11").Value = '  +1.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.350'
$ws.Range("E12").Value = '  +4.54%  '

$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000189'
$ws.Range("E14").Value = '  +10.71%  '

$ws.Range("E15").Value = '  +3.68%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.67'
$ws.Range("E16").Value = '  +3.85%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.932.59'
$ws.Range("E17").Value = '  +4.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.620.53'
$ws.Range("E18").Value = '  +4.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.07'
$ws.Range("E19").Value = '  +7.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '385.16'
$ws.Range("E20").Value = '  +6.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.54'
$ws.Range("E21").Value = '  +6.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.17'
$ws.Range("E22").Value = '  +4.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.01'
$ws.Range("E23").Value = '  +22.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.93'
$ws.Range("E24").Value = '  +4.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.47'
$ws.Range("E25").Value = '  +7.72%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.95'
$ws.Range("E27").Value = '  +12.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.767.14'
$ws.Range("E28").Value = '  +4.76%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0971'
$ws.Range("E30").Value = '  +11.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '549.06'
$ws.Range("E31").Value = '  +8.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.09'
$ws.Range("E32").Value = '  +5.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.34'
$ws.Range("E33").Value = '  +10.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  +3.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '166.03'
$ws.Range("E36").Value = '  +2.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.30'
$ws.Range("E37").Value = '  +4.21%  '

$ws.Range("E38").Value = '  -1.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.14'
$ws.Range("E39").Value = '  +2.71%  '

$ws.Range("E40").Value = '  +7.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.88'
$ws.Range("E41").Value = '  +10.90%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.63'
$ws.Range("E42").Value = '  +13.62%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.05'
$ws.Range("E44").Value = '  +7.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.335'
$ws.Range("E45").Value = '  +5.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.64'
$ws.Range("E46").Value = '  +1.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.16'
$ws.Range("E47").Value = '  +0.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.68'
$ws.Range("E48").Value = '  +3.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.539'
$ws.Range("E49").Value = '  +5.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.70'
$ws.Range("E50").Value = '  +9.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0266'
$ws.Range("E51").Value = '  +7.13%  '
